# Apply the commit change:
#  - TEST_CASES sheet's TC_KIND value (X2) is changed from "SCRIPTED" to "TRAP !!!"
#    (matching the scripting language value already in Y2), which causes the now
#    unused "SCRIPTED" shared string to be dropped from the workbook's string table.
#  - The sheet view's scroll/selection state is updated (topLeftCell K1 -> N1,
#    selection W10 -> X5).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("TEST_CASES")

# Update TC_KIND (column X) value on the data row (row 2) to match the
# TC_SCRIPTING_LANGUAGE value already present, simulating the "incorrect" test type.
$ws.Range("X2").Value = "TRAP !!!"

# Update the saved view state for the sheet: scrolled position and active selection.
$ws.Activate()
$excel.ActiveWindow.ScrollRow = 1
$excel.ActiveWindow.ScrollColumn = 14   # column N is the 14th column
$ws.Range("X5").Select()
